# 11.b.2 indicator sheet — extend the data table with years 2020-2023 and
# convert the 2019 "number of local governments" cell from a text value to
# a real number, matching the upstream author's update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Layout tweaks: narrower name/title columns, taller (wrapped) header row ---
$ws.Columns("A:C").ColumnWidth = 39.33
$ws.Rows.Item(1).RowHeight = 79.5

# --- Extend the year columns: copy the 2019 column's formatting (D3:D6) into
#     the new E:H columns so the new cells inherit the same per-row number
#     formats/styles (general / 0.0% style / general / general) ---
$ws.Range("D3:D6").Copy()
$ws.Range("E3:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 3 — year headers
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# Row 4 — number of local governments (484, unchanged across years); D4 was
# stored as the text "484" and is now a genuine number like the new cells
$ws.Range("D4").Value = 484
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# Row 5 — proportion of local governments (%)
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# Row 6 — number of local governments that adopted/implemented local DRR strategies
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169

# --- Reset selection back to the top-left cell (previously F6) ---
$null = $ws.Range("A1").Select()
